# Weekly refresh of "Hortaliza, Terminal Hortofrutícola Agro Chillán - Berenjena" data.
# Each existing record (rows 2-16) is rotated to a new date/row position (as the
# underlying daily feed advances to the next week), one new record is appended
# as the new row 17, and one row (previously row 8) is replaced by a brand new
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44266
$ws.Cells.Item(2, 11).Value = 9000
$ws.Cells.Item(2, 12).Value = 9500
$ws.Cells.Item(2, 13).Value = 9208
$ws.Cells.Item(2, 15).Value = "Región del Maule"
$ws.Cells.Item(2, 16).Value = 153

# Row 3
$ws.Cells.Item(3, 4).Value = 44218
$ws.Cells.Item(3, 10).Value = 65
$ws.Cells.Item(3, 11).Value = 9000
$ws.Cells.Item(3, 13).Value = 9615
$ws.Cells.Item(3, 16).Value = 160

# Row 4
$ws.Cells.Item(4, 4).Value = 44259
$ws.Cells.Item(4, 10).Value = 70
$ws.Cells.Item(4, 11).Value = 9000
$ws.Cells.Item(4, 12).Value = 9500
$ws.Cells.Item(4, 13).Value = 9214
$ws.Cells.Item(4, 16).Value = 154

# Row 5
$ws.Cells.Item(5, 4).Value = 44224
$ws.Cells.Item(5, 10).Value = 80
$ws.Cells.Item(5, 11).Value = 8500
$ws.Cells.Item(5, 12).Value = 9000
$ws.Cells.Item(5, 13).Value = 8719
$ws.Cells.Item(5, 16).Value = 145

# Row 6
$ws.Cells.Item(6, 4).Value = 44216
$ws.Cells.Item(6, 11).Value = 9500
$ws.Cells.Item(6, 12).Value = 10000
$ws.Cells.Item(6, 13).Value = 9773
$ws.Cells.Item(6, 16).Value = 163

# Row 7
$ws.Cells.Item(7, 4).Value = 44253
$ws.Cells.Item(7, 10).Value = 95
$ws.Cells.Item(7, 11).Value = 9500
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = 9658
$ws.Cells.Item(7, 16).Value = 161

# Row 8
$ws.Cells.Item(8, 4).Value = 44594
$ws.Cells.Item(8, 10).Value = 80
$ws.Cells.Item(8, 11).Value = 12000
$ws.Cells.Item(8, 12).Value = 13000
$ws.Cells.Item(8, 13).Value = 12500
$ws.Cells.Item(8, 16).Value = 208

# Row 9
$ws.Cells.Item(9, 4).Value = 44159
$ws.Cells.Item(9, 10).Value = 35
$ws.Cells.Item(9, 11).Value = 7500
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 7714
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 129

# Row 10
$ws.Cells.Item(10, 4).Value = 44202
$ws.Cells.Item(10, 10).Value = 50
$ws.Cells.Item(10, 11).Value = 8000
$ws.Cells.Item(10, 12).Value = 9000
$ws.Cells.Item(10, 13).Value = 8400
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 140

# Row 11
$ws.Cells.Item(11, 4).Value = 44204
$ws.Cells.Item(11, 10).Value = 45
$ws.Cells.Item(11, 13).Value = 9722
$ws.Cells.Item(11, 16).Value = 162

# Row 12
$ws.Cells.Item(12, 4).Value = 44210
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 9000
$ws.Cells.Item(12, 13).Value = 8417
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 140

# Row 13
$ws.Cells.Item(13, 4).Value = 44271
$ws.Cells.Item(13, 10).Value = 55
$ws.Cells.Item(13, 11).Value = 9000
$ws.Cells.Item(13, 12).Value = 9500
$ws.Cells.Item(13, 13).Value = 9227
$ws.Cells.Item(13, 16).Value = 154

# Row 14
$ws.Cells.Item(14, 4).Value = 44162
$ws.Cells.Item(14, 10).Value = 43
$ws.Cells.Item(14, 11).Value = 8000
$ws.Cells.Item(14, 12).Value = 8500
$ws.Cells.Item(14, 13).Value = 8209
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 137

# Row 15
$ws.Cells.Item(15, 4).Value = 44208
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 8000
$ws.Cells.Item(15, 13).Value = 7350
$ws.Cells.Item(15, 16).Value = 122

# Row 16
$ws.Cells.Item(16, 4).Value = 44264
$ws.Cells.Item(16, 11).Value = 8500
$ws.Cells.Item(16, 12).Value = 9000
$ws.Cells.Item(16, 13).Value = 8709
$ws.Cells.Item(16, 15).Value = "Región del Maule"
$ws.Cells.Item(16, 16).Value = 145

# New row 17 (appended record)
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 44160
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 100112001
$ws.Cells.Item(17, 7).Value = "Berenjena"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 90
$ws.Cells.Item(17, 11).Value = 7500
$ws.Cells.Item(17, 12).Value = 8000
$ws.Cells.Item(17, 13).Value = 7667
$ws.Cells.Item(17, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(17, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value = 128
$ws.Cells.Item(17, 17).Value = 60
$ws.Cells.Item(17, 18).Value = "Hortaliza"
